$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1399
$ws.Range("J19").Value = 1278.625
$ws.Range("L19").Value = 1278.625
$ws.Range("N19").Value = -1628.625
$ws.Range("H32").Value = 9253.77
$ws.Range("J32").Value = 12057.143
$ws.Range("L32").Value = 12057.143
$ws.Range("N32").Value = -12709.143
$ws.Range("H40").Value = 5000.161
$ws.Range("J40").Value = 5354.522
$ws.Range("L40").Value = 5354.522
$ws.Range("N40").Value = -5704.522
$ws.Range("H80").Value = 456.85715
$ws.Range("I80").Value = 378
$ws.Range("K80").Value = 1134
$ws.Range("M80").Value = -136
$ws.Range("H83").Value = 456.85715
$ws.Range("I83").Value = 378
$ws.Range("K83").Value = 3402
$ws.Range("M83").Value = 1590
$ws.Range("H132").Value = 1113.8235
$ws.Range("I132").Value = 927.4
$ws.Range("K132").Value = 2782.2
$ws.Range("M132").Value = -252.1999999999998
$ws.Range("H137").Value = 2975
$ws.Range("I137").Value = 2450
$ws.Range("J137").Value = 3500
$ws.Range("K137").Value = 7350
$ws.Range("L137").Value = 10500
$ws.Range("M137").Value = -4800
$ws.Range("N137").Value = -15600
$ws.Range("H138").Value = 4009.2
$ws.Range("I138").Value = 7546.3335
$ws.Range("J138").Value = 3385
$ws.Range("K138").Value = 22639.0005
$ws.Range("L138").Value = 10155
$ws.Range("M138").Value = -17499.0005
$ws.Range("N138").Value = -20435
$ws.Range("H141").Value = 3584.6667
$ws.Range("I141").Value = 3513.95
$ws.Range("K141").Value = 10541.85
$ws.Range("M141").Value = -5361.849999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5383088.5
$ws.Range("I61").Value = 7581342
$ws.Range("K61").Value = 7581342
$ws.Range("M61").Value = -7581130
$ws.Range("H63").Value = 5177
$ws.Range("J63").Value = 6185
$ws.Range("L63").Value = 6185
$ws.Range("N63").Value = -7557
$ws.Range("H66").Value = 5177
$ws.Range("J66").Value = 6185
$ws.Range("L66").Value = 30925
$ws.Range("N66").Value = -37789
$ws.Range("H74").Value = 4263.241
$ws.Range("I74").Value = 2728.1052
$ws.Range("J74").Value = 7180
$ws.Range("K74").Value = 2728.1052
$ws.Range("L74").Value = 7180
$ws.Range("M74").Value = -1854.1052
$ws.Range("N74").Value = -8928
$ws.Range("H77").Value = 4263.241
$ws.Range("I77").Value = 2728.1052
$ws.Range("J77").Value = 7180
$ws.Range("K77").Value = 13640.526
$ws.Range("L77").Value = 35900
$ws.Range("M77").Value = -9272.526
$ws.Range("N77").Value = -44636
$ws.Range("H136").Value = 5383088.5
$ws.Range("I136").Value = 7581342
$ws.Range("K136").Value = 22744026
$ws.Range("M136").Value = -22741476
$ws.Range("H139").Value = 81342.8
$ws.Range("J139").Value = 81342.8
$ws.Range("L139").Value = 81342.8
$ws.Range("N139").Value = -91622.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4399.6
$ws.Range("I134").Value = 4399.6
$ws.Range("K134").Value = 13198.8
$ws.Range("M134").Value = -10663.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12500
$ws.Range("I58").Value = 7500
$ws.Range("K58").Value = 7500
$ws.Range("M58").Value = -7297
$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52080
$ws.Range("H122").Value = 2938.3447
$ws.Range("I122").Value = 3051.0417
$ws.Range("K122").Value = 9153.125100000001
$ws.Range("M122").Value = -6703.125100000001
$ws.Range("H134").Value = 10516.706
$ws.Range("I134").Value = 7466.3335
$ws.Range("J134").Value = 11170.357
$ws.Range("K134").Value = 22399.0005
$ws.Range("L134").Value = 33511.071
$ws.Range("M134").Value = -19864.0005
$ws.Range("N134").Value = -38581.071
$ws.Range("H135").Value = 75862.94500000001
$ws.Range("J135").Value = 75862.94500000001
$ws.Range("L135").Value = 75862.94500000001
$ws.Range("N135").Value = -86002.94500000001
$ws.Range("H136").Value = 12500
$ws.Range("I136").Value = 7500
$ws.Range("K136").Value = 22500
$ws.Range("M136").Value = -19950

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 475.7619
$ws.Range("I2").Value = 541
$ws.Range("J2").Value = 84.333336
$ws.Range("K2").Value = 3246
$ws.Range("L2").Value = 506.000016
$ws.Range("M2").Value = -3133
$ws.Range("N2").Value = -732.000016
$ws.Range("H5").Value = 1554.6571
$ws.Range("I5").Value = 1431.7059
$ws.Range("J5").Value = 1670.7778
$ws.Range("K5").Value = 4295.1177
$ws.Range("L5").Value = 5012.3334
$ws.Range("M5").Value = -4183.1177
$ws.Range("N5").Value = -5236.3334
$ws.Range("H8").Value = 831.5
$ws.Range("I8").Value = 831.5
$ws.Range("K8").Value = 2494.5
$ws.Range("M8").Value = -2355.5
$ws.Range("H23").Value = 233.15384
$ws.Range("I23").Value = 177.25
$ws.Range("J23").Value = 258
$ws.Range("K23").Value = 531.75
$ws.Range("L23").Value = 774
$ws.Range("M23").Value = -296.75
$ws.Range("N23").Value = -1244
$ws.Range("H68").Value = 3089.3333
$ws.Range("I68").Value = 3691.0908
$ws.Range("J68").Value = 1434.5
$ws.Range("K68").Value = 11073.2724
$ws.Range("L68").Value = 4303.5
$ws.Range("M68").Value = -10262.2724
$ws.Range("N68").Value = -5925.5
$ws.Range("H71").Value = 3089.3333
$ws.Range("I71").Value = 3691.0908
$ws.Range("J71").Value = 1434.5
$ws.Range("K71").Value = 33219.8172
$ws.Range("L71").Value = 12910.5
$ws.Range("M71").Value = -29163.8172
$ws.Range("N71").Value = -21022.5
$ws.Range("H122").Value = 2662.7334
$ws.Range("J122").Value = 2754.2
$ws.Range("L122").Value = 24787.8
$ws.Range("N122").Value = -29687.8
$ws.Range("H126").Value = 14643.333
$ws.Range("I126").Value = 13930
$ws.Range("K126").Value = 41790
$ws.Range("M126").Value = -36850
$ws.Range("H129").Value = 37047820
$ws.Range("I129").Value = 55558316
$ws.Range("K129").Value = 166674948
$ws.Range("M129").Value = -166669948
$ws.Range("H131").Value = 45457170
$ws.Range("I131").Value = 500000000
$ws.Range("J131").Value = 2883
$ws.Range("K131").Value = 1500000000
$ws.Range("L131").Value = 8649
$ws.Range("M131").Value = -1499994960
$ws.Range("N131").Value = -18729
$ws.Range("H135").Value = 1554.6571
$ws.Range("I135").Value = 1431.7059
$ws.Range("J135").Value = 1670.7778
$ws.Range("K135").Value = 12885.3531
$ws.Range("L135").Value = 15037.0002
$ws.Range("M135").Value = -10350.3531
$ws.Range("N135").Value = -20107.0002
$ws.Range("H137").Value = 9872.846
$ws.Range("I137").Value = 2034.8
$ws.Range("K137").Value = 6104.4
$ws.Range("M137").Value = -1004.4
$ws.Range("H140").Value = 11065
$ws.Range("I140").Value = 11065
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 33195
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -28015
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 2750
$ws.Range("I141").Value = 2750
$ws.Range("K141").Value = 8250
$ws.Range("M141").Value = -3070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13573.611
$ws.Range("I70").Value = 10050.571
$ws.Range("K70").Value = 10050.571
$ws.Range("M70").Value = -9780.571
$ws.Range("H73").Value = 13573.611
$ws.Range("I73").Value = 10050.571
$ws.Range("K73").Value = 10050.571
$ws.Range("M73").Value = -9114.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2693.7646
$ws.Range("I46").Value = 1524.625
$ws.Range("K46").Value = 1524.625
$ws.Range("M46").Value = -1336.625
$ws.Range("H136").Value = 4598.5713
$ws.Range("I136").Value = 4450
$ws.Range("K136").Value = 13350
$ws.Range("M136").Value = -10800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 4950
$ws.Range("J26").Value = 4950
$ws.Range("L26").Value = 4950
$ws.Range("N26").Value = -5536
$ws.Range("H41").Value = 22618
$ws.Range("J41").Value = 22618
$ws.Range("L41").Value = 22618
$ws.Range("N41").Value = -23398
$ws.Range("H75").Value = 80001
$ws.Range("I75").Value = 80001
$ws.Range("K75").Value = 80001
$ws.Range("M75").Value = -79065
$ws.Range("H78").Value = 80001
$ws.Range("I78").Value = 80001
$ws.Range("K78").Value = 240003
$ws.Range("M78").Value = -235323
$ws.Range("H107").Value = 594.3077
$ws.Range("I107").Value = 623.9167
$ws.Range("K107").Value = 1871.7501
$ws.Range("M107").Value = 48.24990000000003
$ws.Range("H109").Value = 69996.664
$ws.Range("J109").Value = 69996.664
$ws.Range("L109").Value = 69996.664
$ws.Range("N109").Value = -72770.664
$ws.Range("H122").Value = 4508.826
$ws.Range("J122").Value = 5191
$ws.Range("L122").Value = 15573
$ws.Range("N122").Value = -20473
